# Insert two new price records at the top of the data block (rows 561-562),
# pushing the existing rows 561.. down by two. This mirrors the diff, which
# shows the pre-existing rows 561-603 reappearing unchanged at 563-605, and
# two brand-new rows with fresh data taking the former 561-562 slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before row 561 (old data shifts down by 2).
$ws.Range("A561:R562").EntireRow.Insert()

# --- New row 561 ---
$ws.Cells.Item(561, 1).Value = 9
$ws.Cells.Item(561, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(561, 3).Value = "Metropolitana"
$ws.Cells.Item(561, 4).Value = 45265
$ws.Cells.Item(561, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(561, 5).Value = 13
$ws.Cells.Item(561, 6).Value = 100112017
$ws.Cells.Item(561, 7).Value = "Apio"
$ws.Cells.Item(561, 8).Value = "Americana (o)"
$ws.Cells.Item(561, 9).Value = "Primera"
$ws.Cells.Item(561, 10).Value = 160
$ws.Cells.Item(561, 11).Value = 9000
$ws.Cells.Item(561, 12).Value = 10000
$ws.Cells.Item(561, 13).Value = 9500
$ws.Cells.Item(561, 14).Value = "`$/docena de matas"
$ws.Cells.Item(561, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(561, 16).Value = 1583
$ws.Cells.Item(561, 17).Value = 6
$ws.Cells.Item(561, 18).Value = "Hortaliza"

# --- New row 562 ---
$ws.Cells.Item(562, 1).Value = 9
$ws.Cells.Item(562, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(562, 3).Value = "Metropolitana"
$ws.Cells.Item(562, 4).Value = 45265
$ws.Cells.Item(562, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(562, 5).Value = 13
$ws.Cells.Item(562, 6).Value = 100112017
$ws.Cells.Item(562, 7).Value = "Apio"
$ws.Cells.Item(562, 8).Value = "Americana (o)"
$ws.Cells.Item(562, 9).Value = "Segunda"
$ws.Cells.Item(562, 10).Value = 70
$ws.Cells.Item(562, 11).Value = 8000
$ws.Cells.Item(562, 12).Value = 8000
$ws.Cells.Item(562, 13).Value = 8000
$ws.Cells.Item(562, 14).Value = "`$/docena de matas"
$ws.Cells.Item(562, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(562, 16).Value = 1333
$ws.Cells.Item(562, 17).Value = 6
$ws.Cells.Item(562, 18).Value = "Hortaliza"
